# Natmi following Dr Hou advice
# Recompute the Robo1-Robo1 LR-pair sheet using EC (endothelial cell) data in
# addition to the existing FAPs / sCs / Robo1 clusters, expanding the 3x3 grid
# of sending/target clusters to rows 2-10 and replacing the ligand/receptor
# "expressing cluster" columns (B, C) with the new "ECs" cluster label.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# row 2: FAPs -> FAPs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "ECs"
$ws.Range("C2").Value = "ECs"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = [double]"2"
$ws.Range("F2").Value = [double]"0.6666666666666666"
$ws.Range("G2").Value = [double]"0.1757713333333334"
$ws.Range("H2").Value = [double]"0.5273140000000001"
$ws.Range("I2").Value = [double]"0.009444264308298454"
$ws.Range("J2").Value = [double]"0.009444264308298454"
$ws.Range("K2").Value = [double]"2"
$ws.Range("L2").Value = [double]"0.6666666666666666"
$ws.Range("M2").Value = [double]"0.1757713333333334"
$ws.Range("N2").Value = [double]"0.5273140000000001"
$ws.Range("O2").Value = [double]"0.009444264308298454"
$ws.Range("P2").Value = [double]"0.009444264308298454"
$ws.Range("Q2").Value = [double]"0.03089556162177779"
$ws.Range("R2").Value = [double]"0.2780600545960001"
$ws.Range("S2").Value = [double]"8.919412832500007E-05"
$ws.Range("T2").Value = [double]"8.919412832500007E-05"

# row 3: FAPs -> sCs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "ECs"
$ws.Range("C3").Value = "ECs"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = [double]"2"
$ws.Range("F3").Value = [double]"0.6666666666666666"
$ws.Range("G3").Value = [double]"0.1757713333333334"
$ws.Range("H3").Value = [double]"0.5273140000000001"
$ws.Range("I3").Value = [double]"0.009444264308298454"
$ws.Range("J3").Value = [double]"0.009444264308298454"
$ws.Range("K3").Value = [double]"3"
$ws.Range("L3").Value = [double]"1"
$ws.Range("M3").Value = [double]"16.98312366666667"
$ws.Range("N3").Value = [double]"50.949371"
$ws.Range("O3").Value = [double]"0.9125100529581165"
$ws.Range("P3").Value = [double]"0.9125100529581165"
$ws.Range("Q3").Value = [double]"2.985146291054889"
$ws.Range("R3").Value = [double]"26.866316619494"
$ws.Range("S3").Value = [double]"0.008617986124115871"
$ws.Range("T3").Value = [double]"0.008617986124115871"

# row 4: FAPs -> Robo1
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "ECs"
$ws.Range("C4").Value = "ECs"
$ws.Range("D4").Value = "Robo1"
$ws.Range("E4").Value = [double]"2"
$ws.Range("F4").Value = [double]"0.6666666666666666"
$ws.Range("G4").Value = [double]"0.1757713333333334"
$ws.Range("H4").Value = [double]"0.5273140000000001"
$ws.Range("I4").Value = [double]"0.009444264308298454"
$ws.Range("J4").Value = [double]"0.009444264308298454"
$ws.Range("K4").Value = [double]"3"
$ws.Range("L4").Value = [double]"1"
$ws.Range("M4").Value = [double]"1.452542333333333"
$ws.Range("N4").Value = [double]"4.357627"
$ws.Range("O4").Value = [double]"0.07804568273358503"
$ws.Range("P4").Value = [double]"0.07804568273358505"
$ws.Range("Q4").Value = [double]"0.2553153026531111"
$ws.Range("R4").Value = [double]"2.297837723878"
$ws.Range("S4").Value = [double]"0.000737084055857582"
$ws.Range("T4").Value = [double]"0.0007370840558575822"

# row 5: sCs -> FAPs
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "ECs"
$ws.Range("C5").Value = "ECs"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = [double]"3"
$ws.Range("F5").Value = [double]"1"
$ws.Range("G5").Value = [double]"16.98312366666667"
$ws.Range("H5").Value = [double]"50.949371"
$ws.Range("I5").Value = [double]"0.9125100529581165"
$ws.Range("J5").Value = [double]"0.9125100529581165"
$ws.Range("K5").Value = [double]"2"
$ws.Range("L5").Value = [double]"0.6666666666666666"
$ws.Range("M5").Value = [double]"0.1757713333333334"
$ws.Range("N5").Value = [double]"0.5273140000000001"
$ws.Range("O5").Value = [double]"0.009444264308298454"
$ws.Range("P5").Value = [double]"0.009444264308298454"
$ws.Range("Q5").Value = [double]"2.985146291054889"
$ws.Range("R5").Value = [double]"26.866316619494"
$ws.Range("S5").Value = [double]"0.008617986124115871"
$ws.Range("T5").Value = [double]"0.008617986124115871"

# row 6: sCs -> sCs
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "ECs"
$ws.Range("C6").Value = "ECs"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = [double]"3"
$ws.Range("F6").Value = [double]"1"
$ws.Range("G6").Value = [double]"16.98312366666667"
$ws.Range("H6").Value = [double]"50.949371"
$ws.Range("I6").Value = [double]"0.9125100529581165"
$ws.Range("J6").Value = [double]"0.9125100529581165"
$ws.Range("K6").Value = [double]"3"
$ws.Range("L6").Value = [double]"1"
$ws.Range("M6").Value = [double]"16.98312366666667"
$ws.Range("N6").Value = [double]"50.949371"
$ws.Range("O6").Value = [double]"0.9125100529581165"
$ws.Range("P6").Value = [double]"0.9125100529581165"
$ws.Range("Q6").Value = [double]"288.4264894772935"
$ws.Range("R6").Value = [double]"2595.838405295641"
$ws.Range("S6").Value = [double]"0.8326745967496245"
$ws.Range("T6").Value = [double]"0.8326745967496245"

# row 7: sCs -> Robo1
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "ECs"
$ws.Range("C7").Value = "ECs"
$ws.Range("D7").Value = "Robo1"
$ws.Range("E7").Value = [double]"3"
$ws.Range("F7").Value = [double]"1"
$ws.Range("G7").Value = [double]"16.98312366666667"
$ws.Range("H7").Value = [double]"50.949371"
$ws.Range("I7").Value = [double]"0.9125100529581165"
$ws.Range("J7").Value = [double]"0.9125100529581165"
$ws.Range("K7").Value = [double]"3"
$ws.Range("L7").Value = [double]"1"
$ws.Range("M7").Value = [double]"1.452542333333333"
$ws.Range("N7").Value = [double]"4.357627"
$ws.Range("O7").Value = [double]"0.07804568273358503"
$ws.Range("P7").Value = [double]"0.07804568273358505"
$ws.Range("Q7").Value = [double]"24.66870607806856"
$ws.Range("R7").Value = [double]"222.018354702617"
$ws.Range("S7").Value = [double]"0.07121747008437604"
$ws.Range("T7").Value = [double]"0.07121747008437605"

# row 8: Robo1 -> FAPs
$ws.Range("A8").Value = "Robo1"
$ws.Range("B8").Value = "ECs"
$ws.Range("C8").Value = "ECs"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = [double]"3"
$ws.Range("F8").Value = [double]"1"
$ws.Range("G8").Value = [double]"1.452542333333333"
$ws.Range("H8").Value = [double]"4.357627"
$ws.Range("I8").Value = [double]"0.07804568273358503"
$ws.Range("J8").Value = [double]"0.07804568273358505"
$ws.Range("K8").Value = [double]"2"
$ws.Range("L8").Value = [double]"0.6666666666666666"
$ws.Range("M8").Value = [double]"0.1757713333333334"
$ws.Range("N8").Value = [double]"0.5273140000000001"
$ws.Range("O8").Value = [double]"0.009444264308298454"
$ws.Range("P8").Value = [double]"0.009444264308298454"
$ws.Range("Q8").Value = [double]"0.2553153026531111"
$ws.Range("R8").Value = [double]"2.297837723878"
$ws.Range("S8").Value = [double]"0.000737084055857582"
$ws.Range("T8").Value = [double]"0.0007370840558575822"

# row 9: Robo1 -> sCs
$ws.Range("A9").Value = "Robo1"
$ws.Range("B9").Value = "ECs"
$ws.Range("C9").Value = "ECs"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = [double]"3"
$ws.Range("F9").Value = [double]"1"
$ws.Range("G9").Value = [double]"1.452542333333333"
$ws.Range("H9").Value = [double]"4.357627"
$ws.Range("I9").Value = [double]"0.07804568273358503"
$ws.Range("J9").Value = [double]"0.07804568273358505"
$ws.Range("K9").Value = [double]"3"
$ws.Range("L9").Value = [double]"1"
$ws.Range("M9").Value = [double]"16.98312366666667"
$ws.Range("N9").Value = [double]"50.949371"
$ws.Range("O9").Value = [double]"0.9125100529581165"
$ws.Range("P9").Value = [double]"0.9125100529581165"
$ws.Range("Q9").Value = [double]"24.66870607806856"
$ws.Range("R9").Value = [double]"222.018354702617"
$ws.Range("S9").Value = [double]"0.07121747008437604"
$ws.Range("T9").Value = [double]"0.07121747008437605"

# row 10: Robo1 -> Robo1
$ws.Range("A10").Value = "Robo1"
$ws.Range("B10").Value = "ECs"
$ws.Range("C10").Value = "ECs"
$ws.Range("D10").Value = "Robo1"
$ws.Range("E10").Value = [double]"3"
$ws.Range("F10").Value = [double]"1"
$ws.Range("G10").Value = [double]"1.452542333333333"
$ws.Range("H10").Value = [double]"4.357627"
$ws.Range("I10").Value = [double]"0.07804568273358503"
$ws.Range("J10").Value = [double]"0.07804568273358505"
$ws.Range("K10").Value = [double]"3"
$ws.Range("L10").Value = [double]"1"
$ws.Range("M10").Value = [double]"1.452542333333333"
$ws.Range("N10").Value = [double]"4.357627"
$ws.Range("O10").Value = [double]"0.07804568273358503"
$ws.Range("P10").Value = [double]"0.07804568273358505"
$ws.Range("Q10").Value = [double]"2.109879230125445"
$ws.Range("R10").Value = [double]"18.988913071129"
$ws.Range("S10").Value = [double]"0.006091128593351413"
$ws.Range("T10").Value = [double]"0.006091128593351415"

